# Update LR-pair TPM-derived metrics for Gnai2-Igf1r sheet.
# The underlying change is a refresh of the TPM input values for the
# "ECs" cluster (both as sending cluster ligand-expression and as
# target cluster receptor-expression), which cascades into every
# derived specificity / weight column (G,H,I,J,M,N,O,P,Q,R,S,T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 82.98768099999999
    "H2"  = 248.963043
    "I2"  = 0.4489504115427952
    "J2"  = 0.4489504115427952
    "M2"  = 14.129345
    "N2"  = 42.388035
    "O2"  = 0.3414817166893976
    "P2"  = 0.3414817166893976
    "Q2"  = 1172.561575598945
    "R2"  = 10553.05418039051
    "S2"  = 0.1533083572420453
    "T2"  = 0.1533083572420453

    "G3"  = 82.98768099999999
    "H3"  = 248.963043
    "I3"  = 0.4489504115427952
    "J3"  = 0.4489504115427952
    "O3"  = 0.3553528814026711
    "P3"  = 0.3553528814026711
    "Q3"  = 1220.191635882326
    "R3"  = 10981.72472294094
    "S3"  = 0.1595358223486473
    "T3"  = 0.1595358223486473

    "G4"  = 82.98768099999999
    "H4"  = 248.963043
    "I4"  = 0.4489504115427952
    "J4"  = 0.4489504115427952
    "O4"  = 0.3031654019079313
    "P4"  = 0.3031654019079312
    "Q4"  = 1040.993072116907
    "R4"  = 9368.937649052163
    "S4"  = 0.1361062319521027
    "T4"  = 0.1361062319521026

    "G5"  = 63.14058933333333
    "I5"  = 0.3415807409566563
    "J5"  = 0.3415807409566563
    "M5"  = 14.129345
    "N5"  = 42.388035
    "O5"  = 0.3414817166893976
    "P5"  = 0.3414817166893976
    "Q5"  = 892.1351701939867
    "R5"  = 8029.216531745879
    "S5"  = 0.1166435778099154
    "T5"  = 0.1166435778099154

    "G6"  = 63.14058933333333
    "I6"  = 0.3415807409566563
    "J6"  = 0.3415807409566563
    "O6"  = 0.3553528814026711
    "P6"  = 0.3553528814026711
    "Q6"  = 928.3741642234123
    "S6"  = 0.1213817005306072
    "T6"  = 0.1213817005306072

    "G7"  = 63.14058933333333
    "I7"  = 0.3415807409566563
    "J7"  = 0.3415807409566563
    "O7"  = 0.3031654019079313
    "P7"  = 0.3031654019079312
    "Q7"  = 792.0322061460987
    "R7"  = 7128.289855314888
    "S7"  = 0.1035554626161337
    "T7"  = 0.1035554626161337

    "I8"  = 0.2094688475005485
    "J8"  = 0.2094688475005485
    "M8"  = 14.129345
    "N8"  = 42.388035
    "O8"  = 0.3414817166893976
    "P8"  = 0.3414817166893976
    "Q8"  = 547.0874188979901
    "R8"  = 4923.786770081911
    "S8"  = 0.07152978163743692
    "T8"  = 0.07152978163743692

    "I9"  = 0.2094688475005485
    "J9"  = 0.2094688475005485
    "O9"  = 0.3553528814026711
    "P9"  = 0.3553528814026711
    "Q9"  = 569.3103940360594
    "R9"  = 5123.793546324535
    "S9"  = 0.07443535852341661
    "T9"  = 0.07443535852341661

    "I10" = 0.2094688475005485
    "J10" = 0.2094688475005485
    "O10" = 0.3031654019079313
    "P10" = 0.3031654019079312
    "Q10" = 485.7009002910741
    "R10" = 4371.308102619667
    "S10" = 0.06350370733969495
    "T10" = 0.06350370733969493
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
